$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1869.8158
$ws.Range("J17").Value = 1650.1082
$ws.Range("L17").Value = 4950.3246
$ws.Range("N17").Value = -5286.3246
# Row 64
$ws.Range("H64").Value = 5702.6
$ws.Range("I64").Value = 4750.273
$ws.Range("J64").Value = 6866.5557
$ws.Range("K64").Value = 4750.273
$ws.Range("L64").Value = 6866.5557
$ws.Range("M64").Value = -4502.273
$ws.Range("N64").Value = -7362.5557
# Row 67
$ws.Range("H67").Value = 5702.6
$ws.Range("I67").Value = 4750.273
$ws.Range("J67").Value = 6866.5557
$ws.Range("K67").Value = 4750.273
$ws.Range("L67").Value = 6866.5557
$ws.Range("M67").Value = -3892.273
$ws.Range("N67").Value = -8582.555700000001
# Row 103
$ws.Range("H103").Value = 681.7857
$ws.Range("I103").Value = 486
$ws.Range("J103").Value = 1034.2
$ws.Range("K103").Value = 1458
$ws.Range("L103").Value = 3102.6
$ws.Range("M103").Value = -872
$ws.Range("N103").Value = -4274.6
# Row 112
$ws.Range("H112").Value = 1858.7742
$ws.Range("J112").Value = 1966.5769
$ws.Range("L112").Value = 5899.7307
$ws.Range("N112").Value = -8115.7307
# Row 125
$ws.Range("H125").Value = 1487.7
$ws.Range("J125").Value = 1557.4
$ws.Range("L125").Value = 14016.6
$ws.Range("N125").Value = -18936.6
# Row 132
$ws.Range("H132").Value = 323956.88
$ws.Range("I132").Value = 1422.1
$ws.Range("K132").Value = 4266.299999999999
$ws.Range("M132").Value = -1736.299999999999
# Row 137
$ws.Range("H137").Value = 3999.389
$ws.Range("I137").Value = 3133.8333
$ws.Range("K137").Value = 9401.499899999999
$ws.Range("M137").Value = -6851.499899999999
# Row 138
$ws.Range("H138").Value = 2474.2827
$ws.Range("J138").Value = 2615.3242
$ws.Range("L138").Value = 7845.9726
$ws.Range("N138").Value = -18125.9726

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8258.666999999999
$ws.Range("I32").Value = 5723.7393
$ws.Range("J32").Value = 14089
$ws.Range("K32").Value = 5723.7393
$ws.Range("L32").Value = 14089
$ws.Range("M32").Value = -5436.7393
$ws.Range("N32").Value = -14663
# Row 61
$ws.Range("H61").Value = 4662.8335
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 4662.8335
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 4662.8335
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -5086.8335
# Row 74
$ws.Range("H74").Value = 1199.4
$ws.Range("I74").Value = 999.25
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 999.25
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -125.25
$ws.Range("N74").Value = -3748
# Row 77
$ws.Range("H77").Value = 1199.4
$ws.Range("I77").Value = 999.25
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 4996.25
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -628.25
$ws.Range("N77").Value = -18736
# Row 97
$ws.Range("H97").Value = 1621.9445
$ws.Range("I97").Value = 1698.5
$ws.Range("K97").Value = 1698.5
$ws.Range("M97").Value = -1202.5
# Row 132
$ws.Range("H132").Value = 1577
$ws.Range("I132").Value = 1577
$ws.Range("K132").Value = 4731
$ws.Range("M132").Value = -2201
# Row 136
$ws.Range("H136").Value = 4662.8335
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4662.8335
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 13988.5005
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -19088.5005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3834.5334
$ws.Range("I134").Value = 2621.0476
$ws.Range("K134").Value = 7863.1428
$ws.Range("M134").Value = -5328.1428
# Row 141
$ws.Range("H141").Value = 75499.5
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 33894.5
$ws.Range("J9").Value = 33894.5
$ws.Range("L9").Value = 33894.5
$ws.Range("N9").Value = -34230.5
# Row 132
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1333.3334
$ws.Range("K132").Value = 4000.0002
$ws.Range("M132").Value = -1470.0002
# Row 134
$ws.Range("H134").Value = 9900.454
$ws.Range("I134").Value = 9868.375
$ws.Range("K134").Value = 29605.125
$ws.Range("M134").Value = -27070.125
# Row 141
$ws.Range("H141").Value = 469790.4
$ws.Range("J141").Value = 469790.4
$ws.Range("L141").Value = 469790.4
$ws.Range("N141").Value = -480150.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 94.07407000000001
$ws.Range("J2").Value = 99
$ws.Range("L2").Value = 594
$ws.Range("N2").Value = -820
# Row 37
$ws.Range("H37").Value = 98211.44500000001
$ws.Range("J37").Value = 98211.44500000001
$ws.Range("L37").Value = 294634.335
$ws.Range("N37").Value = -294858.335
# Row 97
$ws.Range("H97").Value = 1097.75
$ws.Range("J97").Value = 1097.5
$ws.Range("L97").Value = 3292.5
$ws.Range("N97").Value = -4284.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Range("H27").Value = 5687.5
$ws.Range("I27").Value = 4250
$ws.Range("K27").Value = 4250
$ws.Range("M27").Value = -4084
# Row 54
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -20780
# Row 70
$ws.Range("H70").Value = 8757
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 9549.833000000001
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 9549.833000000001
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -10089.833
# Row 73
$ws.Range("H73").Value = 8757
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 9549.833000000001
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 9549.833000000001
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -11421.833
# Row 132
$ws.Range("H132").Value = 2557.4285
$ws.Range("I132").Value = 2573.3635
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 7720.0905
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -5190.0905
$ws.Range("N132").Value = -12557

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7192.316
$ws.Range("I7").Value = 8238.429
$ws.Range("J7").Value = 6582.0835
$ws.Range("K7").Value = 8238.429
$ws.Range("L7").Value = 6582.0835
$ws.Range("M7").Value = -8126.429
$ws.Range("N7").Value = -6806.0835
# Row 122
$ws.Range("H122").Value = 4675
$ws.Range("I122").Value = 2400
$ws.Range("K122").Value = 7200
$ws.Range("M122").Value = -4750
# Row 126
$ws.Range("H126").Value = 7192.316
$ws.Range("I126").Value = 8238.429
$ws.Range("J126").Value = 6582.0835
$ws.Range("K126").Value = 24715.287
$ws.Range("L126").Value = 19746.2505
$ws.Range("M126").Value = -22245.287
$ws.Range("N126").Value = -24686.2505
# Row 132
$ws.Range("H132").Value = 3676.5386
$ws.Range("I132").Value = 3779.9
$ws.Range("J132").Value = 3332
$ws.Range("K132").Value = 11339.7
$ws.Range("L132").Value = 9996
$ws.Range("M132").Value = -8809.700000000001
$ws.Range("N132").Value = -15056
# Row 136
$ws.Range("H136").Value = 4393.4
$ws.Range("I136").Value = 4393.4
$ws.Range("K136").Value = 13180.2
$ws.Range("M136").Value = -10630.2
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 29298.5
$ws.Range("J101").Value = 29298.5
$ws.Range("L101").Value = 29298.5
$ws.Range("N101").Value = -35788.5
# Row 122
$ws.Range("H122").Value = 6087.4194
$ws.Range("I122").Value = 5508.885
$ws.Range("K122").Value = 16526.655
$ws.Range("M122").Value = -14076.655
# Row 132
$ws.Range("H132").Value = 3528.5
$ws.Range("I132").Value = 3528.5
$ws.Range("K132").Value = 10585.5
$ws.Range("M132").Value = -8055.5
# Row 136
$ws.Range("H136").Value = 2569.111
$ws.Range("I136").Value = 2569.111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7707.333
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5157.333
$ws.Range("N136").ClearContents()
